# Update column C ("Fitness") values for rows 2-136 on the active sheet.
# The diff groups consecutive rows sharing the same new fitness value;
# each group is written in one shot via a Range assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    @{ Start = 2;   End = 9;   Value = 9166 },
    @{ Start = 10;  End = 12;  Value = 8920 },
    @{ Start = 13;  End = 17;  Value = 8670 },
    @{ Start = 18;  End = 23;  Value = 8356 },
    @{ Start = 24;  End = 43;  Value = 7948 },
    @{ Start = 44;  End = 53;  Value = 7736 },
    @{ Start = 54;  End = 65;  Value = 7345 },
    @{ Start = 66;  End = 85;  Value = 7343 },
    @{ Start = 86;  End = 136; Value = 7310 }
)

foreach ($g in $groups) {
    $rangeAddr = "C$($g.Start):C$($g.End)"
    $ws.Range($rangeAddr).Value = $g.Value
}
